$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K column values (battlecry_mechanic / rerolls fixes)
$ws.Range("K4").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("K13").Value = 0

# Update the active cell selection to K6
$ws.Range("K6").Select()

$wb.Save()
